$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.564.19'
$ws.Cells.Item(2, 5).Value = '  +1.20%  '
$ws.Cells.Item(3, 4).Value = '1.924.78'
$ws.Cells.Item(3, 5).Value = '  +2.98%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  -1.47%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '246.42'
$ws.Cells.Item(5, 5).Value = '  +4.29%  '
$ws.Cells.Item(6, 5).Value = '  -1.30%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4744'
$ws.Cells.Item(7, 5).Value = '  +2.05%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2916'
$ws.Cells.Item(8, 5).Value = '  +4.67%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06790'
$ws.Cells.Item(9, 5).Value = '  +7.05%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '106.37'
$ws.Cells.Item(10, 5).Value = '  +13.83%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '18.39'
$ws.Cells.Item(11, 5).Value = '  +3.45%  '
$ws.Cells.Item(12, 4).Value = '1.912.71'
$ws.Cells.Item(12, 5).Value = '  +1.31%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.07723'
$ws.Cells.Item(13, 5).Value = '  +2.08%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.357'
$ws.Cells.Item(14, 5).Value = '  +8.17%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6723'
$ws.Cells.Item(15, 5).Value = '  +6.38%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '288.80'
$ws.Cells.Item(16, 5).Value = '  +1.95%  '
$ws.Cells.Item(17, 4).Value = '30.615.01'
$ws.Cells.Item(17, 5).Value = '  +1.10%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000007648'
$ws.Cells.Item(18, 5).Value = '  +4.90%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.04'
$ws.Cells.Item(19, 5).Value = '  +3.37%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '1.000'
$ws.Cells.Item(20, 5).Value = '  -0.87%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.501'
$ws.Cells.Item(21, 5).Value = '  +10.29%  '
$ws.Cells.Item(22, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(22, 4).Value = '2.163.45'
$ws.Cells.Item(22, 5).Value = '  -0.92%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.001'
$ws.Cells.Item(23, 5).Value = '  -2.11%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.295'
$ws.Cells.Item(24, 5).Value = '  +5.72%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.415'
$ws.Cells.Item(25, 5).Value = '  +4.22%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '168.48'
$ws.Cells.Item(26, 5).Value = '  +2.42%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.83'
$ws.Cells.Item(27, 5).Value = '  +9.80%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.139'
$ws.Cells.Item(28, 5).Value = '  +12.48%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.1088'
$ws.Cells.Item(29, 5).Value = '  +2.03%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.362'
$ws.Cells.Item(30, 5).Value = '  +0.92%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.183'
$ws.Cells.Item(31, 5).Value = '  +4.85%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.074'
$ws.Cells.Item(32, 5).Value = '  +8.20%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05072'
$ws.Cells.Item(33, 5).Value = '  +3.63%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.7418'
$ws.Cells.Item(34, 5).Value = '  +4.09%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.158'
$ws.Cells.Item(35, 5).Value = '  +4.95%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02083'
$ws.Cells.Item(36, 5).Value = '  +9.82%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.749'
$ws.Cells.Item(37, 5).Value = '  -0.05%  '
$ws.Cells.Item(38, 5).Value = '  +0.02%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.070'
$ws.Cells.Item(39, 5).Value = '  +6.04%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '111.05'
$ws.Cells.Item(40, 5).Value = '  +5.39%  '
$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.8778'
$ws.Cells.Item(41, 5).Value = '  +2.14%  '
$ws.Cells.Item(42, 2).Value = 'TheSandbox'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.4445'
$ws.Cells.Item(42, 5).Value = '  +10.63%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.931'
$ws.Cells.Item(43, 5).Value = '  +7.03%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.000'
$ws.Cells.Item(44, 5).Value = '  -1.35%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '67.60'
$ws.Cells.Item(45, 5).Value = '  +4.71%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '7.297'
$ws.Cells.Item(46, 5).Value = '  +4.37%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '9.322'
$ws.Cells.Item(47, 5).Value = '  +5.57%  '
$ws.Cells.Item(48, 2).Value = 'Algorand'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.1235'
$ws.Cells.Item(48, 5).Value = '  +5.24%  '
$ws.Cells.Item(49, 2).Value = 'Decentraland'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.4145'
$ws.Cells.Item(49, 5).Value = '  +12.50%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '47.25'
$ws.Cells.Item(50, 5).Value = '  +18.86%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '35.18'
$ws.Cells.Item(51, 5).Value = '  +4.90%  '
